$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '250.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("F2").Value = '21-12-2022'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '1'
$ws.Range("G2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").Value = '21-12-2022'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '1'
$ws.Range("G3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.455'
$ws.Range("D4").Style = "Normal"
$ws.Range("F4").Value = '21-12-2022'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '1'
$ws.Range("G4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05756'
$ws.Range("D5").Style = "Normal"
$ws.Range("F5").Value = '21-12-2022'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '1'
$ws.Range("G5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.421'
$ws.Range("D6").Style = "Normal"
$ws.Range("F6").Value = '21-12-2022'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '1'
$ws.Range("G6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.339'
$ws.Range("D7").Style = "Normal"
$ws.Range("F7").Value = '21-12-2022'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '1'
$ws.Range("G7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8150'
$ws.Range("D8").Style = "Normal"
$ws.Range("F8").Value = '21-12-2022'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '1'
$ws.Range("G8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9317'
$ws.Range("D9").Style = "Normal"
$ws.Range("F9").Value = '21-12-2022'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '1'
$ws.Range("G9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1435'
$ws.Range("D10").Style = "Normal"
$ws.Range("F10").Value = '21-12-2022'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '1'
$ws.Range("G10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07512'
$ws.Range("D11").Style = "Normal"
$ws.Range("F11").Value = '21-12-2022'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '1'
$ws.Range("G11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03120'
$ws.Range("D12").Style = "Normal"
$ws.Range("F12").Value = '21-12-2022'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '1'
$ws.Range("G12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03105'
$ws.Range("D13").Style = "Normal"
$ws.Range("F13").Value = '21-12-2022'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '1'
$ws.Range("G13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09363'
$ws.Range("D14").Style = "Normal"
$ws.Range("F14").Value = '21-12-2022'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '1'
$ws.Range("G14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.717'
$ws.Range("D15").Style = "Normal"
$ws.Range("F15").Value = '21-12-2022'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '1'
$ws.Range("G15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001573'
$ws.Range("D16").Style = "Normal"
$ws.Range("F16").Value = '21-12-2022'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '1'
$ws.Range("G16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04755'
$ws.Range("D17").Style = "Normal"
$ws.Range("F17").Value = '21-12-2022'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '1'
$ws.Range("G17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005800'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17OneONE'
$ws.Range("F18").Value = '21-12-2022'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '1'
$ws.Range("G18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006433'
$ws.Range("D19").Style = "Normal"
$ws.Range("F19").Value = '21-12-2022'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '1'
$ws.Range("G19").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005004'
$ws.Range("D20").Style = "Normal"
$ws.Range("F20").Value = '21-12-2022'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '1'
$ws.Range("G20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001029'
$ws.Range("D21").Style = "Normal"
$ws.Range("F21").Value = '21-12-2022'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '1'
$ws.Range("G21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001501'
$ws.Range("D22").Style = "Normal"
$ws.Range("F22").Value = '21-12-2022'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '1'
$ws.Range("G22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.700'
$ws.Range("D23").Style = "Normal"
$ws.Range("F23").Value = '21-12-2022'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '1'
$ws.Range("G23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.152'
$ws.Range("D24").Style = "Normal"
$ws.Range("F24").Value = '21-12-2022'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '1'
$ws.Range("G24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3303'
$ws.Range("D25").Style = "Normal"
$ws.Range("F25").Value = '21-12-2022'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '1'
$ws.Range("G25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1310'
$ws.Range("D26").Style = "Normal"
$ws.Range("F26").Value = '21-12-2022'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '1'
$ws.Range("G26").Style = "Normal"

# Row 27
$ws.Range("F27").Value = '21-12-2022'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '1'
$ws.Range("G27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0003005'
$ws.Range("D28").Style = "Normal"
$ws.Range("F28").Value = '21-12-2022'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '1'
$ws.Range("G28").Style = "Normal"

# Row 29
$ws.Range("F29").Value = '21-12-2022'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '1'
$ws.Range("G29").Style = "Normal"

# Row 30
$ws.Range("F30").Value = '21-12-2022'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '1'
$ws.Range("G30").Style = "Normal"

# Row 31
$ws.Range("F31").Value = '21-12-2022'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '1'
$ws.Range("G31").Style = "Normal"

# Row 32
$ws.Range("F32").Value = '21-12-2022'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '1'
$ws.Range("G32").Style = "Normal"

# Row 33
$ws.Range("F33").Value = '21-12-2022'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '1'
$ws.Range("G33").Style = "Normal"

# Row 34
$ws.Range("F34").Value = '21-12-2022'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '1'
$ws.Range("G34").Style = "Normal"

# Row 35
$ws.Range("F35").Value = '21-12-2022'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '1'
$ws.Range("G35").Style = "Normal"

# Row 36
$ws.Range("F36").Value = '21-12-2022'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '1'
$ws.Range("G36").Style = "Normal"

# Row 37
$ws.Range("F37").Value = '21-12-2022'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '1'
$ws.Range("G37").Style = "Normal"

# Row 38
$ws.Range("F38").Value = '21-12-2022'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '1'
$ws.Range("G38").Style = "Normal"

# Row 39
$ws.Range("F39").Value = '21-12-2022'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '1'
$ws.Range("G39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04014'
$ws.Range("D40").Style = "Normal"
$ws.Range("F40").Value = '21-12-2022'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '1'
$ws.Range("G40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1067'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("F41").Value = '21-12-2022'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '1'
$ws.Range("G41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002712'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("F42").Value = '21-12-2022'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '1'
$ws.Range("G42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002948'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$ws.Range("F43").Value = '21-12-2022'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '1'
$ws.Range("G43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008004'
$ws.Range("D44").Style = "Normal"
$ws.Range("F44").Value = '21-12-2022'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '1'
$ws.Range("G44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005901'
$ws.Range("D45").Style = "Normal"
$ws.Range("F45").Value = '21-12-2022'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '1'
$ws.Range("G45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000751'
$ws.Range("D46").Style = "Normal"
$ws.Range("F46").Value = '21-12-2022'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '1'
$ws.Range("G46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5009'
$ws.Range("D47").Style = "Normal"
$ws.Range("F47").Value = '21-12-2022'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '1'
$ws.Range("G47").Style = "Normal"

# Row 48
$ws.Range("F48").Value = '21-12-2022'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '1'
$ws.Range("G48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002104'
$ws.Range("D49").Style = "Normal"
$ws.Range("F49").Value = '21-12-2022'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '1'
$ws.Range("G49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01012'
$ws.Range("D50").Style = "Normal"
$ws.Range("F50").Value = '21-12-2022'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '1'
$ws.Range("G50").Style = "Normal"

# Row 51
$ws.Range("F51").Value = '21-12-2022'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '1'
$ws.Range("G51").Style = "Normal"
